$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Candidate details ---
$ws.Range("C3").Value = "Mayur"
$ws.Range("C4").Value = "Chavan"
$ws.Range("C5").Value = "Support Engineer"
$ws.Range("C7").Value = "Deepak Yadav"
$ws.Range("C8").Value = "Michelin India"

# --- Joining date ---
$ws.Range("C10").Value = "8/19/2024"

# --- Address ---
$ws.Range("C11").Value = "202,Wing A,MadhuMalti Apartment"
$ws.Range("D11").Value = "Shubhash Nagar,Belwali,Badlapur West,421503"

# --- Contact details ---
$ws.Range("C12").Value = 9405210295

# --- Number of dependents ---
$ws.Range("C13").Value = 0

# --- Dependent travel details ---
$ws.Range("C18").Value = "8/19/2024"
$ws.Range("C19").Value = "Cab "
$ws.Range("C20").Value = "Badlpaur (Mumbai)"
$ws.Range("C21").Value = "Kharadi (Pune)"
$ws.Range("C22").Value = "8/18/2024"
$ws.Range("C23").Value = "Y"
$ws.Range("C24").Value = "Yes(Only Vehicle required , TVS Jupiter scooty)"

# Row 24 was manually resized by the author after the text was entered.
$ws.Rows.Item(24).RowHeight = 52.5

# --- Requested date (apply explicit date number format, matching d-mmm-yy) ---
$ws.Range("C27").Value = "8/7/2024"
$ws.Range("C27").NumberFormat = "d-mmm-yy"

# Reflect the author's final selection in the sheet view.
$ws.Range("C28").Select()
